# Stand Up Meeting Week 8 — fill in Guillermo Toloza Guzman's rows (16-18)
# with his answers, mirroring the pattern already used for the other
# team members higher up in the sheet. Also adds the new shared string
# "Planeación de reunión" (used in F16) and updates the sheet's
# selection state to match where the author left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# Row 16 - "¿Qué se hizo ayer?"
$ws.Range("C16").Value = "Nada "
$ws.Range("D16").Value = "Nada"
$ws.Range("E16").Value = "Nada"
$ws.Range("F16").Value = "Planeación de reunión"
$ws.Range("G16").Value = "Se logró el objetivo y se encontró el diagrama que faltaba"

# Row 17 - "¿Qué se hará hoy?"
$ws.Range("C17").Value = "Nada "
$ws.Range("D17").Value = "Nada"
$ws.Range("E17").Value = "Planeamos reunión para mañana"
$ws.Range("F17").Value = "Reunirnos, revisar diagramas y Mariana comenzará con las tablas"
$ws.Range("G17").Value = "Reunión para revisar y opinar el trabajo de Mariana"

# Row 18 - "¿Qué cosas se oponen?"
$ws.Range("C18").Value = "Nada "
$ws.Range("D18").Value = "No hubo clase "
$ws.Range("E18").Value = "Nada"
$ws.Range("F18").Value = "Nada"
$ws.Range("G18").Value = "Parcial mañana"

# Reflect the author's final scroll/selection position in the sheet view.
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("G19").Select()
